# Daily attendance processing - 2025-11-14 06:32:13
# Swap the order of the two comma-separated "modified by" entries in
# column G ("Modified By") for the specific rows that changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows (in column G) whose two comma-separated values need to swap order.
$rows = @(3,6,7,10,11,12,13,14,15,17,18,19,20,21,22,24,26,29,32,33,36,37,38,39,40,41,43,44,45,46,47,48,50,52,55,58,59,62,63,64,65,66,67,69,70,71,72,73,74,76,78,83,84,85,86,87,90,92,93,94,96,99,101,109,110,111,112,113,116,118,119,120,122,125,127,135,136,137,138,139,142,144,145,146,148,151,153)

foreach ($r in $rows) {
    $cell = $ws.Range("G$r")
    $val = [string]$cell.Text

    $commaIndex = $val.IndexOf(",")
    if ($commaIndex -ge 0) {
        $first = $val.Substring(0, $commaIndex).Trim()
        $second = $val.Substring($commaIndex + 1).Trim()
        $cell.Value = "$second, $first"
    }
}
